# The forecast window rolled forward by 2 days (the sheet always holds the most
# recent 169 hourly rows): row 2 used to start at 30.11.2025 14:00 and now
# starts at 02.12.2025 13:00. Every row keeps its row-local "Interval" cycle
# (1..24, day rolls over after 24), and the "Prediction" (xgb forecast output)
# was also refreshed for this run - most quiet night-time hours stay at the
# floor value of 0.011 but the values around each day's solar peak change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow   = 2
$lastRow    = 170
$rowCount   = $lastRow - $firstRow + 1

$startSerial = 45993   # 02.12.2025 (Excel serial date, epoch 1899-12-30)
$startHour   = 13

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

# Updated "Prediction" values (row number -> new value) - every other row
# keeps whatever value it already had before this edit.
$overrides = @{
    4   = 0.768
    5   = 0.73
    6   = 0.249
    7   = 0.019
    21  = 0.011
    22  = 0.018
    23  = 0.37
    24  = 1.228
    25  = 1.602
    26  = 1.53
    27  = 1.769
    28  = 1.811
    29  = 1.292
    30  = 0.29
    31  = 0.016
    45  = 0.011
    46  = 0.018
    47  = 0.337
    48  = 1.042
    49  = 1.619
    50  = 2.237
    51  = 2.217
    52  = 1.616
    53  = 1.158
    54  = 0.4
    55  = 0.022
    69  = 0.011
    70  = 0.011
    71  = 0.152
    72  = 0.618
    73  = 1.364
    74  = 1.189
    75  = 1.173
    76  = 0.752
    77  = 0.486
    78  = 0.117
    79  = 0.014
    94  = 0.011
    95  = 0.08
    96  = 0.241
    97  = 0.343
    98  = 0.453
    99  = 0.65
    100 = 0.519
    101 = 0.347
    102 = 0.092
    118 = 0.011
    119 = 0.08
    120 = 0.299
    121 = 0.5
    122 = 0.926
    123 = 0.835
    124 = 0.538
    125 = 0.307
    126 = 0.092
    142 = 0.011
    143 = 0.137
    144 = 0.561
    145 = 0.922
    146 = 1.158
    147 = 1.144
    148 = 0.711
    149 = 0.396
    150 = 0.118
    151 = 0.014
    166 = 0.011
    167 = 0.132
    168 = 0.364
    169 = 0.777
    170 = 1.239
}

$data = New-Object 'object[,]' $rowCount,4

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $firstRow + $i

    $totalHour = $startHour + $i
    $dayOffset = [Math]::Floor(($totalHour - 1) / 24)
    $hour      = $totalHour - ($dayOffset * 24)
    $serial    = $startSerial + $dayOffset

    $prediction = $ws.Cells.Item($row, 3).Value
    if ($overrides.ContainsKey($row)) {
        $prediction = $overrides[$row]
    }

    $dateText = $epoch.AddDays($serial).ToString("dd.MM.yyyy")
    $lookup   = "$dateText$hour"

    $data[$i,0] = $serial
    $data[$i,1] = $hour
    $data[$i,2] = $prediction
    $data[$i,3] = $lookup
}

$ws.Range("A$firstRow`:D$lastRow").Value = $data
